# Add 2022-02-24 data: update nombre_aides (col C) and montant_total (col E)
# for the rows whose underlying source figures changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;   C = 100802; E = 327280426 }
    @{ Row = 3;   C = 249276; E = 1036206026 }
    @{ Row = 5;   C = 39478;  E = 361437789 }
    @{ Row = 46;  C = 10921;  E = 66194111 }
    @{ Row = 53;  C = 141656; E = 589988411 }
    @{ Row = 55;  C = 23188;  E = 187779348 }
    @{ Row = 57;  C = 3706;   E = 137946554 }
    @{ Row = 63;  C = 14102;  E = 35461428 }
    @{ Row = 64;  C = 5056;   E = 19410015 }
    @{ Row = 79;  C = 116575; E = 447286129 }
    @{ Row = 81;  C = 17424;  E = 133408450 }
    @{ Row = 91;  C = 150987; E = 480953671 }
    @{ Row = 92;  C = 408712; E = 1590457459 }
    @{ Row = 93;  C = 209259; E = 1304015728 }
    @{ Row = 94;  C = 94004;  E = 911508273 }
    @{ Row = 96;  C = 17146;  E = 780965357 }
    @{ Row = 104; C = 135154; E = 271719428 }
    @{ Row = 106; C = 18119;  E = 40714435 }
    @{ Row = 114; C = 3715;   E = 8916220 }
    @{ Row = 115; C = 11471;  E = 32308090 }
    @{ Row = 116; C = 4426;   E = 19536136 }
    @{ Row = 118; C = 908;    E = 10462697 }
    @{ Row = 122; C = 8329;   E = 12560350 }
    @{ Row = 131; C = 75574;  E = 307077678 }
    @{ Row = 166; C = 35926;  E = 210544120 }
    @{ Row = 174; C = 226047; E = 900134076 }
    @{ Row = 175; C = 80751;  E = 485295848 }
    @{ Row = 184; C = 68727;  E = 134116344 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
